# Refresh cryptos price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D/E store plain text (prices use '.' as both thousands
# and decimal separator, e.g. "35.560.23"; volumes are padded percent
# strings). Prefixing with a literal apostrophe forces Excel to keep the
# assigned value as text instead of silently re-parsing it as a number
# (which would e.g. turn "95.30" into 95.3 or "11.38" into a numeric cell).

# Row 2
$ws.Range('D2').Value = "'" + '35.560.23'
$ws.Range('E2').Value = "'" + '  +3.34%  '

# Row 3
$ws.Range('D3').Value = "'" + '1.845.09'
$ws.Range('E3').Value = "'" + '  +2.66%  '

# Row 4
$ws.Range('E4').Value = "'" + '  +0.27%  '

# Row 5
$ws.Range('D5').Value = "'" + '231.89'
$ws.Range('E5').Value = "'" + '  +3.43%  '

# Row 6
$ws.Range('D6').Value = "'" + '0.624'
$ws.Range('E6').Value = "'" + '  +3.47%  '

# Row 7
$ws.Range('E7').Value = "'" + '  +0.23%  '

# Row 8
$ws.Range('D8').Value = "'" + '44.46'
$ws.Range('E8').Value = "'" + '  +14.07%  '

# Row 9
$ws.Range('E9').Value = "'" + '  +8.69%  '

# Row 10
$ws.Range('D10').Value = "'" + '0.0699'
$ws.Range('E10').Value = "'" + '  +4.86%  '

# Row 11
$ws.Range('E11').Value = "'" + '  +2.34%  '

# Row 12
$ws.Range('D12').Value = "'" + '2.111.29'
$ws.Range('E12').Value = "'" + '  +2.64%  '

# Row 13
$ws.Range('B13').Value = "'" + 'Chainlink'
$ws.Range('C13').Value = "'" + 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = "'" + '11.38'
$ws.Range('E13').Value = "'" + '  +5.19%  '

# Row 14
$ws.Range('B14').Value = "'" + 'WrappedEther'
$ws.Range('C14').Value = "'" + 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = "'" + '1.842.88'
$ws.Range('E14').Value = "'" + '  +2.74%  '

# Row 15
$ws.Range('D15').Value = "'" + '0.675'
$ws.Range('E15').Value = "'" + '  +7.77%  '

# Row 16
$ws.Range('D16').Value = "'" + '4.74'
$ws.Range('E16').Value = "'" + '  +9.35%  '

# Row 17
$ws.Range('D17').Value = "'" + '35.555.12'
$ws.Range('E17').Value = "'" + '  +3.29%  '

# Row 18
$ws.Range('D18').Value = "'" + '70.62'
$ws.Range('E18').Value = "'" + '  +4.13%  '

# Row 19
$ws.Range('D19').Value = "'" + '0.0₃0805'
$ws.Range('E19').Value = "'" + '  +5.50%  '

# Row 20
$ws.Range('D20').Value = "'" + '244.36'
$ws.Range('E20').Value = "'" + '  +2.32%  '

# Row 21
$ws.Range('D21').Value = "'" + '12.09'
$ws.Range('E21').Value = "'" + '  +10.03%  '

# Row 22
$ws.Range('D22').Value = "'" + '4.65'
$ws.Range('E22').Value = "'" + '  +14.54%  '

# Row 23
$ws.Range('E23').Value = "'" + '  +0.33%  '

# Row 24
$ws.Range('D24').Value = "'" + '2.23'
$ws.Range('E24').Value = "'" + '  +3.07%  '

# Row 25
$ws.Range('D25').Value = "'" + '171.06'
$ws.Range('E25').Value = "'" + '  -0.19%  '

# Row 26
$ws.Range('D26').Value = "'" + '7.99'
$ws.Range('E26').Value = "'" + '  +5.18%  '

# Row 27
$ws.Range('D27').Value = "'" + '17.84'
$ws.Range('E27').Value = "'" + '  +2.38%  '

# Row 28
$ws.Range('D28').Value = "'" + '0.123'
$ws.Range('E28').Value = "'" + '  +1.55%  '

# Row 29
$ws.Range('E29').Value = "'" + '  +29.22%  '

# Row 30
$ws.Range('E30').Value = "'" + '  +0.29%  '

# Row 31
$ws.Range('D31').Value = "'" + '3.341.88'
$ws.Range('E31').Value = "'" + '  +37.54%  '

# Row 32
$ws.Range('D32').Value = "'" + '0.0556'
$ws.Range('E32').Value = "'" + '  +8.59%  '

# Row 33
$ws.Range('D33').Value = "'" + '4.15'
$ws.Range('E33').Value = "'" + '  +9.21%  '

# Row 34
$ws.Range('E34').Value = "'" + '  +6.29%  '

# Row 35
$ws.Range('D35').Value = "'" + '1.84'
$ws.Range('E35').Value = "'" + '  +2.66%  '

# Row 36
$ws.Range('D36').Value = "'" + '95.30'
$ws.Range('E36').Value = "'" + '  +16.92%  '

# Row 37
$ws.Range('B37').Value = "'" + 'ImmutableX'
$ws.Range('C37').Value = "'" + 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = "'" + '0.698'
$ws.Range('E37').Value = "'" + '  +9.90%  '

# Row 38
$ws.Range('B38').Value = "'" + 'TrustWalletToken'
$ws.Range('C38').Value = "'" + 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = "'" + '1.15'
$ws.Range('E38').Value = "'" + '  +9.75%  '

# Row 39
$ws.Range('D39').Value = "'" + '1.345.18'
$ws.Range('E39').Value = "'" + '  +2.97%  '

# Row 40
$ws.Range('D40').Value = "'" + '0.0197'
$ws.Range('E40').Value = "'" + '  +6.83%  '

# Row 41
$ws.Range('E41').Value = "'" + '  +8.25%  '

# Row 42
$ws.Range('E42').Value = "'" + '  +8.83%  '

# Row 43
$ws.Range('D43').Value = "'" + '15.23'
$ws.Range('E43').Value = "'" + '  +7.89%  '

# Row 44
$ws.Range('E44').Value = "'" + '  +3.54%  '

# Row 45
$ws.Range('D45').Value = "'" + '2.45'
$ws.Range('E45').Value = "'" + '  +0.42%  '

# Row 46
$ws.Range('E46').Value = "'" + '  +0.46%  '

# Row 47
$ws.Range('D47').Value = "'" + '6.25'
$ws.Range('E47').Value = "'" + '  +9.94%  '

# Row 48
$ws.Range('D48').Value = "'" + '0.0516'
$ws.Range('E48').Value = "'" + '  +0.07%  '

# Row 49
$ws.Range('D49').Value = "'" + '2.018.67'
$ws.Range('E49').Value = "'" + '  +3.02%  '

# Row 50
$ws.Range('E50').Value = "'" + '  +0.33%  '

# Row 51
$ws.Range('D51').Value = "'" + '102.45'
$ws.Range('E51').Value = "'" + '  +1.05%  '
